$wb = $excel.ActiveWorkbook

# Duplicate the "Erscheinungsverlauf" sheet to create the new
# "Anmerkungen zum Erscheinungsverlauf" sheet, inserted right after it.
$src = $wb.Worksheets.Item("Erscheinungsverlauf")
$src.Copy($null, $src)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
# Excel sheet names are capped at 31 characters; the full title
# "Anmerkungen zum Erscheinungsverlauf" is truncated to that limit.
$ws.Name = "Anmerkungen zum Erscheinungsver"

# Header block (rows 1-6) describing the field.
$ws.Range("B1").Value = "Anmerkungen zum Erscheinungsverlauf"
$ws.Range("B2").Value = "Anzeige der Anmerkungen zum Erscheinungsverlauf"
$ws.Range("B3").Value = "Numbering Peculiarities Note"
$ws.Range("B4").Value = "Anmerkungen zum Erscheinungsverlauf"
$ws.Range("B5").Value = "Numbering peculiarities"
$ws.Range("B6").Value = "Ticket #133"

# Update the hyperlink on the ticket reference to point at ticket #133.
$ws.Hyperlinks.Item(1).TextToDisplay = "Ticket #133"
$ws.Hyperlinks.Item(1).Address = "http://redmine.thulb.uni-jena.de/issues/133"

# Data-field summary row.
$ws.Range("A9").Value = "515 `$a"
$ws.Range("B9").Value = "Anmerkungen zur Zählung von fortlaufenden Ressourcen; ist in Pica 4225 nicht wiederholbar"

# Existing example row now documents the 2008 non-appearance.
$ws.Range("A13").Value = "515 `$a"
$ws.Range("C13").Value = "2008 nicht ersch"
$ws.Range("C13").Style = $ws.Range("A1").Style
$ws.Rows.Item(13).RowHeight = 14.15

# New example row for the irregular-appearance note.
$ws.Range("A14").Value = "515 `$a"
$ws.Range("A14").Style = $ws.Range("A13").Style
$ws.Range("B14").Value = 502081112
$ws.Range("B14").Style = $ws.Range("B13").Style
$ws.Range("C14").Value = "Ersch. unregelmäßig"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("D14").Style = $ws.Range("D13").Style
$ws.Range("E14").Style = $ws.Range("E13").Style
$ws.Range("F14").Style = $ws.Range("F13").Style
$ws.Rows.Item(14).RowHeight = 14.15

# Make the new sheet the active one, as in the authored workbook.
$ws.Activate()
$ws.Select()
